# Add a "Lien_Fichier" column (full file path) before the existing "Lien" column.
# This mirrors the FilesDIR generator being re-run with an extra "path + filename" output column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; this shifts the existing "Lien" column (D) to E, along with its data.
$ws.Columns("D:D").Insert()

# Header for the newly inserted column.
$ws.Range("D1").Value = "Lien_Fichier"

$base = "F:\testBis\dl\Nouveau dossier\_Jeux"

# Rows 2-12: files directly under the base folder
$lien = $base
for ($r = 2; $r -le 12; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 13-16: files under subfolder "avantApres"
$lien = $base + "\avantApres"
for ($r = 13; $r -le 16; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 17-25: files under subfolder "bed"
$lien = $base + "\bed"
for ($r = 17; $r -le 25; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 26-27: files under subfolder "bikiniStreet"
$lien = $base + "\bikiniStreet"
for ($r = 26; $r -le 27; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 28-31: files under subfolder "douche"
$lien = $base + "\douche"
for ($r = 28; $r -le 31; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 32-32: files under subfolder "doucheHabille"
$lien = $base + "\doucheHabille"
for ($r = 32; $r -le 32; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 33-33: files under subfolder "gene"
$lien = $base + "\gene"
for ($r = 33; $r -le 33; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 34-35: files under subfolder "innocent"
$lien = $base + "\innocent"
for ($r = 34; $r -le 35; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 36-43: files under subfolder "irresistible"
$lien = $base + "\irresistible"
for ($r = 36; $r -le 43; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 44-44: files under subfolder "mainSeins"
$lien = $base + "\mainSeins"
for ($r = 44; $r -le 44; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 45-45: files under subfolder "mousse"
$lien = $base + "\mousse"
for ($r = 45; $r -le 45; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 46-46: files under subfolder "nuePro"
$lien = $base + "\nuePro"
for ($r = 46; $r -le 46; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 47-53: files under subfolder "oops"
$lien = $base + "\oops"
for ($r = 47; $r -le 53; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 54-55: files under subfolder "relaxed"
$lien = $base + "\relaxed"
for ($r = 54; $r -le 55; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 56-60: files under subfolder "sexyNoNude"
$lien = $base + "\sexyNoNude"
for ($r = 56; $r -le 60; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 61-66: files under subfolder "sortirDouche"
$lien = $base + "\sortirDouche"
for ($r = 61; $r -le 66; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}

# Rows 67-67: files under subfolder "towel"
$lien = $base + "\towel"
for ($r = 67; $r -le 67; $r++) {
    $fichier = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 5).Value = $lien
    $ws.Cells.Item($r, 4).Value = $lien + "\" + $fichier
}
